$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "IClientBalance-20240730-092431-"

$ws.Cells.Item(2, 7).Value = 45503
$ws.Cells.Item(3, 7).Value = 45503
$ws.Cells.Item(4, 7).Value = 45503
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 7505.95
$ws.Cells.Item(5, 7).Value = 45503
$ws.Cells.Item(6, 7).Value = 45503
$ws.Cells.Item(7, 7).Value = 45503
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 3297.74
$ws.Cells.Item(8, 7).Value = 45503
$ws.Cells.Item(9, 7).Value = 45503
$ws.Cells.Item(10, 7).Value = 45503
$ws.Cells.Item(11, 7).Value = 45503
$ws.Cells.Item(12, 7).Value = 45503
$ws.Cells.Item(13, 7).Value = 45503
$ws.Cells.Item(14, 7).Value = 45503
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 6214.89
$ws.Cells.Item(15, 7).Value = 45503
$ws.Cells.Item(16, 7).Value = 45503
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 4427.76
$ws.Cells.Item(17, 7).Value = 45503
$ws.Cells.Item(18, 7).Value = 45503
$ws.Cells.Item(19, 7).Value = 45503
$ws.Cells.Item(20, 7).Value = 45503
$ws.Cells.Item(21, 7).Value = 45503
$ws.Cells.Item(22, 7).Value = 45503
$ws.Cells.Item(23, 7).Value = 45503
$ws.Cells.Item(24, 7).Value = 45503
$ws.Cells.Item(25, 7).Value = 45503
$ws.Cells.Item(26, 7).Value = 45503
$ws.Cells.Item(27, 7).Value = 45503
$ws.Cells.Item(28, 7).Value = 45503
$ws.Cells.Item(29, 7).Value = 45503
$ws.Cells.Item(30, 7).Value = 45503
$ws.Cells.Item(31, 7).Value = 45503
$ws.Cells.Item(32, 7).Value = 45503
$ws.Cells.Item(33, 7).Value = 45503
$ws.Cells.Item(34, 7).Value = 45503
$ws.Cells.Item(35, 7).Value = 45503
$ws.Cells.Item(36, 5).Value = 503.3
$ws.Cells.Item(36, 7).Value = 45503
$ws.Cells.Item(36, 8).Value = 503.3
$ws.Cells.Item(37, 7).Value = 45503
$ws.Cells.Item(38, 7).Value = 45503
$ws.Cells.Item(39, 7).Value = 45503
$ws.Cells.Item(40, 7).Value = 45503
$ws.Cells.Item(41, 7).Value = 45503
$ws.Cells.Item(42, 7).Value = 45503
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 7692.23
$ws.Cells.Item(43, 7).Value = 45503
$ws.Cells.Item(44, 7).Value = 45503
$ws.Cells.Item(45, 7).Value = 45503
$ws.Cells.Item(46, 7).Value = 45503
$ws.Cells.Item(47, 7).Value = 45503
$ws.Cells.Item(48, 7).Value = 45503
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 2134.1999999999998
$ws.Cells.Item(49, 7).Value = 45503
$ws.Cells.Item(50, 7).Value = 45503
$ws.Cells.Item(51, 7).Value = 45503
$ws.Cells.Item(52, 7).Value = 45503
$ws.Cells.Item(53, 7).Value = 45503
$ws.Cells.Item(54, 7).Value = 45503
$ws.Cells.Item(55, 7).Value = 45503
$ws.Cells.Item(56, 7).Value = 45503
$ws.Cells.Item(57, 7).Value = 45503
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 590.91
$ws.Cells.Item(58, 7).Value = 45503
$ws.Cells.Item(59, 7).Value = 45503
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 6818.84
$ws.Cells.Item(60, 7).Value = 45503
$ws.Cells.Item(61, 7).Value = 45503
$ws.Cells.Item(62, 7).Value = 45503
$ws.Cells.Item(63, 7).Value = 45503
$ws.Cells.Item(64, 7).Value = 45503
$ws.Cells.Item(65, 7).Value = 45503
$ws.Cells.Item(66, 7).Value = 45503
$ws.Cells.Item(67, 7).Value = 45503
$ws.Cells.Item(68, 7).Value = 45503
$ws.Cells.Item(69, 7).Value = 45503
$ws.Cells.Item(70, 7).Value = 45503
$ws.Cells.Item(71, 7).Value = 45503
$ws.Cells.Item(72, 7).Value = 45503
$ws.Cells.Item(73, 7).Value = 45503
$ws.Cells.Item(74, 7).Value = 45503
$ws.Cells.Item(75, 7).Value = 45503
$ws.Cells.Item(76, 7).Value = 45503
$ws.Cells.Item(77, 7).Value = 45503
$ws.Cells.Item(78, 7).Value = 45503
$ws.Cells.Item(79, 7).Value = 45503
$ws.Cells.Item(80, 7).Value = 45503
$ws.Cells.Item(81, 7).Value = 45503
$ws.Cells.Item(82, 7).Value = 45503
$ws.Cells.Item(83, 7).Value = 45503
$ws.Cells.Item(84, 7).Value = 45503
$ws.Cells.Item(85, 7).Value = 45503
$ws.Cells.Item(86, 7).Value = 45503
$ws.Cells.Item(87, 7).Value = 45503
$ws.Cells.Item(88, 7).Value = 45503
$ws.Cells.Item(89, 7).Value = 45503
$ws.Cells.Item(90, 7).Value = 45503
$ws.Cells.Item(91, 7).Value = 45503
$ws.Cells.Item(92, 7).Value = 45503
$ws.Cells.Item(93, 7).Value = 45503
$ws.Cells.Item(94, 7).Value = 45503
$ws.Cells.Item(95, 7).Value = 45503
$ws.Cells.Item(96, 7).Value = 45503
$ws.Cells.Item(97, 7).Value = 45503
$ws.Cells.Item(98, 7).Value = 45503
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 5703.15
$ws.Cells.Item(99, 7).Value = 45503
$ws.Cells.Item(100, 7).Value = 45503
$ws.Cells.Item(101, 7).Value = 45503
$ws.Cells.Item(102, 7).Value = 45503
$ws.Cells.Item(103, 7).Value = 45503
$ws.Cells.Item(104, 4).Value = -739.22
$ws.Cells.Item(104, 5).Value = 19655.009999999998
$ws.Cells.Item(104, 7).Value = 45503
$ws.Cells.Item(104, 8).Value = 18915.79
$ws.Cells.Item(105, 7).Value = 45503
$ws.Cells.Item(106, 7).Value = 45503
$ws.Cells.Item(107, 7).Value = 45503
$ws.Cells.Item(108, 4).Value = -34357.980000000003
$ws.Cells.Item(108, 5).Value = 18028.86
$ws.Cells.Item(108, 7).Value = 45503
$ws.Cells.Item(108, 8).Value = -16329.12
$ws.Cells.Item(109, 7).Value = 45503
$ws.Cells.Item(110, 7).Value = 45503
$ws.Cells.Item(111, 7).Value = 45503
$ws.Cells.Item(112, 7).Value = 45503
$ws.Cells.Item(113, 7).Value = 45503
$ws.Cells.Item(114, 7).Value = 45503
$ws.Cells.Item(115, 7).Value = 45503
$ws.Cells.Item(116, 7).Value = 45503
$ws.Cells.Item(117, 7).Value = 45503
$ws.Cells.Item(118, 7).Value = 45503
$ws.Cells.Item(119, 7).Value = 45503
$ws.Cells.Item(120, 7).Value = 45503
$ws.Cells.Item(121, 7).Value = 45503
$ws.Cells.Item(122, 7).Value = 45503
$ws.Cells.Item(123, 7).Value = 45503
$ws.Cells.Item(124, 7).Value = 45503
$ws.Cells.Item(125, 7).Value = 45503
$ws.Cells.Item(126, 7).Value = 45503
$ws.Cells.Item(127, 7).Value = 45503
$ws.Cells.Item(128, 7).Value = 45503
$ws.Cells.Item(129, 7).Value = 45503
$ws.Cells.Item(130, 7).Value = 45503
$ws.Cells.Item(131, 7).Value = 45503
$ws.Cells.Item(132, 4).Value = 0
$ws.Cells.Item(132, 5).Value = 2501.04
$ws.Cells.Item(132, 7).Value = 45503
$ws.Cells.Item(133, 7).Value = 45503
$ws.Cells.Item(134, 7).Value = 45503
$ws.Cells.Item(135, 7).Value = 45503
$ws.Cells.Item(136, 7).Value = 45503
$ws.Cells.Item(137, 7).Value = 45503
$ws.Cells.Item(138, 7).Value = 45503
$ws.Cells.Item(139, 7).Value = 45503
$ws.Cells.Item(140, 7).Value = 45503
$ws.Cells.Item(141, 7).Value = 45503
$ws.Cells.Item(142, 7).Value = 45503
$ws.Cells.Item(143, 4).Value = 0
$ws.Cells.Item(143, 5).Value = 20298.62
$ws.Cells.Item(143, 7).Value = 45503
$ws.Cells.Item(144, 7).Value = 45503
$ws.Cells.Item(145, 7).Value = 45503
$ws.Cells.Item(146, 7).Value = 45503
$ws.Cells.Item(147, 7).Value = 45503
$ws.Cells.Item(148, 7).Value = 45503
$ws.Cells.Item(149, 7).Value = 45503
$ws.Cells.Item(150, 7).Value = 45503
$ws.Cells.Item(151, 7).Value = 45503
$ws.Cells.Item(152, 7).Value = 45503
$ws.Cells.Item(153, 7).Value = 45503
$ws.Cells.Item(154, 7).Value = 45503
$ws.Cells.Item(155, 7).Value = 45503
$ws.Cells.Item(156, 7).Value = 45503
$ws.Cells.Item(157, 7).Value = 45503
$ws.Cells.Item(158, 4).Value = 0
$ws.Cells.Item(158, 5).Value = 1222.83
$ws.Cells.Item(158, 7).Value = 45503
$ws.Cells.Item(159, 7).Value = 45503
$ws.Cells.Item(160, 7).Value = 45503
$ws.Cells.Item(161, 7).Value = 45503
$ws.Cells.Item(162, 7).Value = 45503
$ws.Cells.Item(163, 7).Value = 45503
$ws.Cells.Item(164, 7).Value = 45503
$ws.Cells.Item(165, 7).Value = 45503
$ws.Cells.Item(166, 7).Value = 45503
$ws.Cells.Item(167, 7).Value = 45503
$ws.Cells.Item(168, 7).Value = 45503
$ws.Cells.Item(169, 7).Value = 45503
$ws.Cells.Item(170, 7).Value = 45503
$ws.Cells.Item(171, 7).Value = 45503
$ws.Cells.Item(172, 4).Value = 0
$ws.Cells.Item(172, 5).Value = 1846.51
$ws.Cells.Item(172, 7).Value = 45503
$ws.Cells.Item(173, 4).Value = -455.51
$ws.Cells.Item(173, 5).Value = 13517.34
$ws.Cells.Item(173, 7).Value = 45503
$ws.Cells.Item(173, 8).Value = 13061.83
$ws.Cells.Item(174, 7).Value = 45503
$ws.Cells.Item(175, 7).Value = 45503
$ws.Cells.Item(176, 7).Value = 45503
$ws.Cells.Item(177, 7).Value = 45503
$ws.Cells.Item(178, 7).Value = 45503
$ws.Cells.Item(179, 7).Value = 45503
$ws.Cells.Item(180, 7).Value = 45503
$ws.Cells.Item(181, 7).Value = 45503
$ws.Cells.Item(182, 7).Value = 45503
$ws.Cells.Item(183, 7).Value = 45503
$ws.Cells.Item(184, 7).Value = 45503
$ws.Cells.Item(185, 7).Value = 45503
$ws.Cells.Item(186, 7).Value = 45503
$ws.Cells.Item(187, 7).Value = 45503
$ws.Cells.Item(188, 7).Value = 45503
$ws.Cells.Item(189, 7).Value = 45503
$ws.Cells.Item(190, 7).Value = 45503
$ws.Cells.Item(191, 7).Value = 45503
$ws.Cells.Item(192, 7).Value = 45503
$ws.Cells.Item(193, 7).Value = 45503
$ws.Cells.Item(194, 7).Value = 45503
$ws.Cells.Item(195, 7).Value = 45503
$ws.Cells.Item(196, 7).Value = 45503
$ws.Cells.Item(197, 7).Value = 45503
$ws.Cells.Item(198, 7).Value = 45503
$ws.Cells.Item(199, 7).Value = 45503
$ws.Cells.Item(200, 7).Value = 45503
$ws.Cells.Item(201, 7).Value = 45503
$ws.Cells.Item(202, 7).Value = 45503
$ws.Cells.Item(203, 7).Value = 45503
$ws.Cells.Item(204, 7).Value = 45503
$ws.Cells.Item(205, 7).Value = 45503
$ws.Cells.Item(206, 7).Value = 45503
$ws.Cells.Item(207, 7).Value = 45503
$ws.Cells.Item(208, 7).Value = 45503
$ws.Cells.Item(209, 7).Value = 45503
$ws.Cells.Item(210, 7).Value = 45503
$ws.Cells.Item(211, 7).Value = 45503
$ws.Cells.Item(212, 7).Value = 45503
$ws.Cells.Item(213, 7).Value = 45503
$ws.Cells.Item(214, 7).Value = 45503
$ws.Cells.Item(215, 7).Value = 45503
$ws.Cells.Item(216, 7).Value = 45503
$ws.Cells.Item(217, 7).Value = 45503
$ws.Cells.Item(218, 7).Value = 45503
$ws.Cells.Item(219, 7).Value = 45503
$ws.Cells.Item(220, 7).Value = 45503
$ws.Cells.Item(221, 7).Value = 45503
$ws.Cells.Item(222, 7).Value = 45503
$ws.Cells.Item(223, 7).Value = 45503
$ws.Cells.Item(224, 7).Value = 45503
$ws.Cells.Item(225, 7).Value = 45503
$ws.Cells.Item(226, 7).Value = 45503
$ws.Cells.Item(227, 7).Value = 45503
$ws.Cells.Item(228, 7).Value = 45503
$ws.Cells.Item(229, 7).Value = 45503
$ws.Cells.Item(230, 7).Value = 45503
$ws.Cells.Item(231, 7).Value = 45503
$ws.Cells.Item(232, 7).Value = 45503
$ws.Cells.Item(233, 7).Value = 45503
$ws.Cells.Item(234, 7).Value = 45503
$ws.Cells.Item(235, 4).Value = 0
$ws.Cells.Item(235, 5).Value = 5199
$ws.Cells.Item(235, 7).Value = 45503
$ws.Cells.Item(236, 7).Value = 45503
$ws.Cells.Item(237, 7).Value = 45503
$ws.Cells.Item(238, 7).Value = 45503
$ws.Cells.Item(239, 7).Value = 45503
$ws.Cells.Item(240, 7).Value = 45503
$ws.Cells.Item(241, 7).Value = 45503
$ws.Cells.Item(242, 7).Value = 45503
$ws.Cells.Item(243, 7).Value = 45503
$ws.Cells.Item(244, 7).Value = 45503
$ws.Cells.Item(245, 7).Value = 45503
$ws.Cells.Item(246, 7).Value = 45503
$ws.Cells.Item(247, 7).Value = 45503
$ws.Cells.Item(248, 7).Value = 45503
$ws.Cells.Item(249, 4).Value = 0
$ws.Cells.Item(249, 5).Value = 5100.8100000000004
$ws.Cells.Item(249, 7).Value = 45503
$ws.Cells.Item(250, 7).Value = 45503
$ws.Cells.Item(251, 7).Value = 45503
$ws.Cells.Item(252, 7).Value = 45503
$ws.Cells.Item(253, 7).Value = 45503
$ws.Cells.Item(254, 7).Value = 45503
$ws.Cells.Item(255, 7).Value = 45503
$ws.Cells.Item(256, 7).Value = 45503
$ws.Cells.Item(257, 7).Value = 45503
$ws.Cells.Item(258, 7).Value = 45503
$ws.Cells.Item(259, 7).Value = 45503
$ws.Cells.Item(260, 7).Value = 45503
$ws.Cells.Item(261, 7).Value = 45503
$ws.Cells.Item(262, 7).Value = 45503
$ws.Cells.Item(263, 7).Value = 45503
$ws.Cells.Item(264, 4).Value = -14579.41
$ws.Cells.Item(264, 5).Value = 15639.47
$ws.Cells.Item(264, 7).Value = 45503
$ws.Cells.Item(264, 8).Value = 1060.06
$ws.Cells.Item(265, 4).Value = 0
$ws.Cells.Item(265, 5).Value = 9294.44
$ws.Cells.Item(265, 7).Value = 45503
$ws.Cells.Item(266, 7).Value = 45503
$ws.Cells.Item(267, 7).Value = 45503
$ws.Cells.Item(268, 7).Value = 45503
$ws.Cells.Item(269, 7).Value = 45503
$ws.Cells.Item(270, 4).Value = 0
$ws.Cells.Item(270, 5).Value = 5760.35
$ws.Cells.Item(270, 7).Value = 45503
$ws.Cells.Item(271, 4).Value = 0
$ws.Cells.Item(271, 5).Value = 7952.12
$ws.Cells.Item(271, 7).Value = 45503
$ws.Cells.Item(272, 7).Value = 45503
$ws.Cells.Item(273, 4).Value = 0
$ws.Cells.Item(273, 5).Value = 5499.66
$ws.Cells.Item(273, 7).Value = 45503
$ws.Cells.Item(274, 7).Value = 45503

# Move active selection back to A1 (default)
$ws.Range("A1").Select()
